# "Changes for the AIT field"
#
# The author renamed the sanity-test fixture from "Automated Sanity" to
# "Automated Sanity2" everywhere it is used across the workbook (the plain
# label cells, plus the longer "Automated Sanity Child2" / "Automated
# Sanity Reseller2" / confirmation-message / "Working as admin ..." strings
# that embed it as a substring). This was done as a single workbook-wide
# Find & Replace ("Replace All", matching on part of the cell, not only
# whole-cell matches), which is why every descendant phrase containing
# "Automated Sanity" picked up the "2" as well.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # LookAt:=xlPart (0) so substrings inside longer phrases (e.g.
    # "Automated Sanity Child2") are updated too, not just exact matches.
    $null = $ws.Cells.Replace("Automated Sanity", "Automated Sanity2", 0, 1, $false, $false, $false, $false)
}

# The active sheet/cell at the time of the edit: ConfigAccType was the
# selected tab, with the cursor left on C1 (the just-edited "Automated
# Sanity2" cell) afterwards.
$active = $wb.Worksheets.Item("ConfigAccType")
$active.Activate()
$active.Range("C1").Select()
